# Generate Report for Handback
# Update the handback timestamps for the c6c61c5a-...-39e8865daba8.md row
# (row 3) on each worksheet, reflecting a refreshed handback/xliff
# generation run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for the
# c6c61c5a-1579-4e80-85c7-39e8865daba8.md row.
$wsOverview.Range("G3").Value = "2016-08-28 02:48:10"

# zh-cn sheet: "Correspond Handoff Datetime" (H) and
# "Correspond Handback DateTime" (K) columns for the same row.
$wsZhCn.Range("H3").Value = "2016-08-28 02:48:06"
$wsZhCn.Range("K3").Value = "2016-08-28 02:48:23"

# de-de sheet: "Correspond Handoff Datetime" (H) and
# "Correspond Handback DateTime" (K) columns for the same row.
$wsDeDe.Range("H3").Value = "2016-08-28 02:48:10"
$wsDeDe.Range("K3").Value = "2016-08-28 02:48:29"
